$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# "add d13C unc nosams": the carbonate_d13C_unc_permil column (U) for the
# NOSAMS-run samples in rows 2-7 held the "NA" placeholder text; fill in the
# reported uncertainty value (0.1 permil) as a number instead.
$ws.Range("U2").Value = 0.1
$ws.Range("U3").Value = 0.1
$ws.Range("U4").Value = 0.1
$ws.Range("U5").Value = 0.1
$ws.Range("U6").Value = 0.1
$ws.Range("U7").Value = 0.1

# Reflect where the user ended up looking/selecting after typing the values.
$ws.Activate()
$ws.Range("U2:U7").Select()
